$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (players reshuffled; Nicolas Claxton replaced by Scotty Pippen Jr.)
$data = @(
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Scotty Pippen Jr.", "PG,SG", "Memphis Grizzlies"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
